$d = $word.ActiveDocument

# Locate the "sta_latitude varchar(40) not null," paragraph (the "stand" table's
# latitude column) rather than the similarly-worded "centroid_latitude" one
# elsewhere in the document, so we search for the precise leading text.
$latIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "sta_latitude varchar*") {
        $latIndex = $i
        break
    }
}

if ($latIndex -eq 0) {
    throw "Could not find the 'sta_latitude varchar(...)' paragraph"
}

$lonIndex = $latIndex + 1
$lonText = $d.Paragraphs.Item($lonIndex).Range.Text
if ($lonText -notlike "sta_longitude varchar*") {
    throw "Paragraph after sta_latitude was not sta_longitude as expected"
}

# Step 1: rename "latitude" -> "coordinates" within the first paragraph only,
# so the spell-check markers around "sta_coordinates" stay correctly paired.
$latRange = $d.Paragraphs.Item($latIndex).Range
$latRange.Find.Execute("latitude", $true, $false, $false, $false, $false, $true, 1, $false, "coordinates", 2) | Out-Null

# Step 2: the two columns "sta_latitude varchar(40) not null," and
# "sta_longitude varchar(40) not null, " collapse into a single
# "sta_coordinates point not null, " column, with the now-empty second
# paragraph left behind (blank line) instead of being removed. Do this as one
# Find/Replace spanning both paragraphs (using ^p for the paragraph mark) so
# that the paragraph break itself is preserved cleanly, without leaving any
# stray empty runs or unmatched spell/grammar-check markers behind.
$spanRange = $d.Range($d.Paragraphs.Item($latIndex).Range.Start, $d.Paragraphs.Item($lonIndex).Range.End)
$spanRange.Find.Execute(" varchar(40) not null,^psta_longitude varchar(40) not null, ", $true, $false, $false, $false, $false, $true, 1, $false, " point not null, ^p", 2) | Out-Null
